# CHORE: Update the Documentserver Docker details
#
# The re-converted copy of this bill that came back from the Documentserver
# pipeline no longer carries the legacy SharePoint "Custom XML Data" parts
# (customXml/item1-3.xml + their customXml/itemProps1-3.xml companions).
# Strip that metadata here the same way Word's own "Custom XML Data" cleanup
# does: walk ActiveDocument.CustomXMLParts and delete every part.

$d = $word.ActiveDocument

$parts = $d.CustomXMLParts
$count = $parts.Count

# Walk backwards so deleting an item never invalidates the index of the
# item(s) we still have to visit.
for ($i = $count; $i -ge 1; $i--) {
    $part = $parts.Item($i)
    $part.Delete()
}

# Defensive second pass in case Delete() renumbers instead of compacting, or
# the host exposes a live/re-queryable collection after each removal.
while ($d.CustomXMLParts.Count -gt 0) {
    $d.CustomXMLParts.Item(1).Delete()
}
